$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records ("Fruta / hortaliza, semanal") need to be
# inserted for Betarraga - Terminal Hortofrutícola Agro Chillán, right
# before the existing row 326. Inserting whole rows there shifts every
# row from the old 326 down through the old 415 two positions down (new
# rows 328-417), which grows the used range from A1:R415 to A1:R417.
$ws.Rows("326:327").Insert()

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# New row 326 - "Primera" quality record dated 2022-08-12 (serial 44785)
$ws.Cells.Item(326, 1).Value = 7
$ws.Cells.Item(326, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(326, 3).Value = "Ñuble"
$ws.Cells.Item(326, 4).Value = 44785
$ws.Cells.Item(326, 4).NumberFormat = $dateFormat
$ws.Cells.Item(326, 5).Value = 16
$ws.Cells.Item(326, 6).Value = 100114014
$ws.Cells.Item(326, 7).Value = "Betarraga"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 300
$ws.Cells.Item(326, 11).Value = 700
$ws.Cells.Item(326, 12).Value = 800
$ws.Cells.Item(326, 13).Value = 750
$ws.Cells.Item(326, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(326, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(326, 16).Value = 150
$ws.Cells.Item(326, 17).Value = 5
$ws.Cells.Item(326, 18).Value = "Hortaliza"

# New row 327 - "Segunda" quality record, same date as row 326
$ws.Cells.Item(327, 1).Value = 7
$ws.Cells.Item(327, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(327, 3).Value = "Ñuble"
$ws.Cells.Item(327, 4).Value = 44785
$ws.Cells.Item(327, 4).NumberFormat = $dateFormat
$ws.Cells.Item(327, 5).Value = 16
$ws.Cells.Item(327, 6).Value = 100114014
$ws.Cells.Item(327, 7).Value = "Betarraga"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Segunda"
$ws.Cells.Item(327, 10).Value = 200
$ws.Cells.Item(327, 11).Value = 600
$ws.Cells.Item(327, 12).Value = 600
$ws.Cells.Item(327, 13).Value = 600
$ws.Cells.Item(327, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(327, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(327, 16).Value = 120
$ws.Cells.Item(327, 17).Value = 5
$ws.Cells.Item(327, 18).Value = "Hortaliza"
